# Replace the text "price" column (column B, rows 2-21) — previously stored
# as shared-string labels like " $   5.000" — with plain numeric values, per
# the new reference price list. Column A (product names) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# row number -> new numeric price
$prices = [ordered]@{
    2  = 5000
    3  = 40000
    4  = 30000
    5  = 30000
    6  = 30000
    7  = 30000
    8  = 1000
    9  = 3000
    10 = 10000
    11 = 17000
    12 = 400
    13 = 10000
    14 = 21000
    15 = 3500
    16 = 6500
    17 = 3000
    18 = 4000
    19 = 2300
    20 = 28000
    21 = 50000
}

foreach ($row in $prices.Keys) {
    $ws.Cells.Item($row, 2).Value = $prices[$row]
}

# Mirror the cursor/selection left behind after the edit.
$ws.Range("B22").Select()
